$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange/UpDown values for the last existing row (row 10)
$ws.Range("X10").Value = -3.4100040000000149
$ws.Range("Y10").Value = "Down"

# Append a brand new trade row (row 11), copying formatting from row 10 above
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("S10").Copy($ws.Range("S11"))
$ws.Range("T10").Copy($ws.Range("T11"))

$ws.Range("A11").Value = 42654.886678240742
$ws.Range("B11").Value = -1
$ws.Range("C11").Value = "Neutral"
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = "Random"
$ws.Range("Q11").Value = 39.313912976930268
$ws.Range("R11").Value = 1.8
$ws.Range("S11").Value = 0.086400000000000005
$ws.Range("T11").Value = -0.0115
$ws.Range("U11").Value = 5.85
$ws.Range("V11").Value = "N/A"
$ws.Range("W11").Value = 0
